$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed as a number by Excel;
# force them to stay as text so formatting like trailing zeros is preserved,
# then restore the default (Normal) style so no stray formatting remains.
$textCells = @("D4", "D5", "D6", "D14", "D19", "D20", "D21", "D24", "D29", "D30", "D32", "D33", "D35", "D39", "D41", "D42", "D43", "D45", "D46", "D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D4").Value = "1.00"
$ws.Range("D5").Value = "592.46"
$ws.Range("D6").Value = "143.66"
$ws.Range("D14").Value = "27.40"
$ws.Range("D19").Value = "11.43"
$ws.Range("D20").Value = "340.26"
$ws.Range("D21").Value = "4.38"
$ws.Range("D24").Value = "67.23"
$ws.Range("D29").Value = "539.41"
$ws.Range("D30").Value = "8.42"
$ws.Range("D32").Value = "1.99"
$ws.Range("D33").Value = "1.81"
$ws.Range("D35").Value = "174.79"
$ws.Range("D39").Value = "19.08"
$ws.Range("D41").Value = "172.27"
$ws.Range("D42").Value = "0.998"
$ws.Range("D43").Value = "40.36"
$ws.Range("D45").Value = "22.21"
$ws.Range("D46").Value = "0.0564"
$ws.Range("D50").Value = "18.65"

foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}

# Remaining cell updates (safe as plain text assignments).
$ws.Range("D2").Value = "63.473.62"
$ws.Range("D3").Value = "2.648.96"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("E5").Value = "  +1.78%  "
$ws.Range("E6").Value = "  -0.87%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  -1.11%  "
$ws.Range("D9").Value = "2.648.44"
$ws.Range("E9").Value = "  +2.54%  "
$ws.Range("E10").Value = "  +0.07%  "
$ws.Range("E11").Value = "  +0.81%  "
$ws.Range("E12").Value = "  +0.99%  "
$ws.Range("E13").Value = "  +0.48%  "
$ws.Range("E14").Value = "  +1.62%  "
$ws.Range("D15").Value = "3.125.63"
$ws.Range("E15").Value = "  +2.57%  "
$ws.Range("D16").Value = "63.404.75"
$ws.Range("E16").Value = "  +1.17%  "
$ws.Range("E17").Value = "  -0.46%  "
$ws.Range("D18").Value = "2.649.94"
$ws.Range("E18").Value = "  +2.23%  "
$ws.Range("E19").Value = "  +1.73%  "
$ws.Range("E20").Value = "  +0.07%  "
$ws.Range("E21").Value = "  +0.59%  "
$ws.Range("E22").Value = "  +1.26%  "
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("E24").Value = "  -0.06%  "
$ws.Range("E25").Value = "  +5.31%  "
$ws.Range("E26").Value = "  -0.47%  "
$ws.Range("E27").Value = "  +2.47%  "
$ws.Range("E28").Value = "  +0.27%  "
$ws.Range("E29").Value = "  +17.60%  "
$ws.Range("E30").Value = "  +2.28%  "
$ws.Range("E31").Value = "  -1.65%  "
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("E32").Value = "  +3.12%  "
$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("E33").Value = "  +11.76%  "
$ws.Range("E34").Value = "  -0.66%  "
$ws.Range("E35").Value = "  -1.12%  "
$ws.Range("E36").Value = "  +9.08%  "
$ws.Range("E37").Value = "  +0.05%  "
$ws.Range("E38").Value = "  +0.21%  "
$ws.Range("E39").Value = "  +0.73%  "
$ws.Range("E40").Value = "  +7.15%  "
$ws.Range("E41").Value = "  +9.23%  "
$ws.Range("E42").Value = "  -0.06%  "
$ws.Range("E43").Value = "  +2.58%  "
$ws.Range("E44").Value = "  +0.36%  "
$ws.Range("E45").Value = "  +4.97%  "
$ws.Range("E46").Value = "  +5.07%  "
$ws.Range("E47").Value = "  +0.74%  "
$ws.Range("E48").Value = "  -0.57%  "
$ws.Range("E49").Value = "  +1.54%  "
$ws.Range("E50").Value = "  +2.66%  "
$ws.Range("E51").Value = "  -0.67%  "
